#
# Adds the final slide ("More Things to Try") to the ML talk deck.
#
$p = $ppt.ActivePresentation

# New slide, inserted at the end, using the "Title and Content" layout
# (same CustomLayout index used by every other content slide in this deck).
$s = $p.Slides.Add($p.Slides.Count + 1, 2)

# --- Title -----------------------------------------------------------
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "More Things to Try"

# --- Body --------------------------------------------------------------
# Build every paragraph first (all at the default outline level) so each
# one gets its own run, then promote/demote the specific paragraphs that
# need a different indent level afterwards - doing it in this order keeps
# newly-appended paragraphs from inheriting an already-changed indent.
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.Text = "See if you can figure out what the rest of the perceptron code is doing"
[void]$body.InsertAfter("`rI’ve commented most of it")
[void]$body.InsertAfter("`rAlter the code to work with a different dataset")
[void]$body.InsertAfter("`rYou will probably have to change some numbers at least")
[void]$body.InsertAfter("`rPractice more things in R")
[void]$body.InsertAfter("`rhttps://www.statmethods.net/r-tutorial/index.html")
[void]$body.InsertAfter("`rhttps://www.cyclismo.org/tutorial/R/")

$body.Paragraphs(2).IndentLevel = 2
$body.Paragraphs(4).IndentLevel = 2
$body.Paragraphs(6).IndentLevel = 2
$body.Paragraphs(7).IndentLevel = 2

# --- Hyperlinks ----------------------------------------------------------
# Paragraph 6: "https://www.statmethods.net/r-tutorial/index.html" - split
# the same way PowerPoint splits auto-detected hyperlink text from text
# typed immediately afterwards, but both runs point at the same link.
$para6 = $body.Paragraphs(6)
$link1a = $para6.Characters(1, 39)
$link1a.ActionSettings.Item(1).Hyperlink.Address = "https://www.statmethods.net/r-tutorial/index.html"
$link1b = $para6.Characters(40, 10)
$link1b.ActionSettings.Item(1).Hyperlink.Address = "https://www.statmethods.net/r-tutorial/index.html"

# Paragraph 7: "https://www.cyclismo.org/tutorial/R/"
$para7 = $body.Paragraphs(7)
$link2a = $para7.Characters(1, 35)
$link2a.ActionSettings.Item(1).Hyperlink.Address = "https://www.cyclismo.org/tutorial/R/"
$link2b = $para7.Characters(36, 1)
$link2b.ActionSettings.Item(1).Hyperlink.Address = "https://www.cyclismo.org/tutorial/R/"
